$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format (bold/centered/bordered style used by the other
# header cells in row 1) onto the new H1 header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column, defaulting to 0 for the existing rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
